$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 17214
$ws.Range("E2").Value = 793
$ws.Range("F2").Value = 793
$ws.Range("G2").Value = 601
$ws.Range("H2").Value = 479
$ws.Range("I2").Value = 475
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 38465
$ws.Range("L2").Value = 27483
$ws.Range("M2").Value = 10981
$ws.Range("N2").Value = 10578
$ws.Range("O2").Value = 403
$ws.Range("P2").Value = 470
$ws.Range("Q2").Value = 1671
$ws.Range("R2").Value = -1895
$ws.Range("S2").Value = 25
$ws.Range("T2").Value = 1307
$ws.Range("U2").Value = 364
$ws.Range("V2").Value = 13170
$ws.Range("W2").Value = 4.6
$ws.Range("X2").Value = 2.78
$ws.Range("AA2").Value = 250.27
$ws.Range("AB2").Value = 2134.61
$ws.Range("AC2").Value = 3054
$ws.Range("AE2").Value = 22591
$ws.Range("AF2").Value = 1.62
$ws.Range("AH2").Value = 1.09
$ws.Range("AI2").Value = 39.4
$ws.Range("AJ2").Value = 46957120
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("AD2").ClearContents()

# Row 3
$ws.Range("D3").Value = 52992
$ws.Range("E3").Value = 2656
$ws.Range("F3").Value = 2656
$ws.Range("G3").Value = 2012
$ws.Range("H3").Value = 1296
$ws.Range("I3").Value = 1258
$ws.Range("J3").Value = 38
$ws.Range("K3").Value = 41598
$ws.Range("L3").Value = 27941
$ws.Range("M3").Value = 13658
$ws.Range("N3").Value = 13219
$ws.Range("O3").Value = 438
$ws.Range("P3").Value = 470
$ws.Range("Q3").Value = 2176
$ws.Range("R3").Value = -2486
$ws.Range("S3").Value = 711
$ws.Range("T3").Value = 2700
$ws.Range("U3").Value = -523
$ws.Range("V3").Value = 12373
$ws.Range("W3").Value = 5.01
$ws.Range("X3").Value = 2.44
$ws.Range("Y3").Value = 10.57
$ws.Range("Z3").Value = 3.24
$ws.Range("AA3").Value = 204.58
$ws.Range("AB3").Value = 2328.59
$ws.Range("AC3").Value = 2679
$ws.Range("AE3").Value = 28232
$ws.Range("AF3").Value = 1.17
$ws.Range("AG3").Value = 960
$ws.Range("AH3").Value = 2.9
$ws.Range("AI3").Value = 35.73
$ws.Range("AJ3").Value = 46957120
$ws.Range("AD3").ClearContents()

# Row 4
$ws.Range("D4").Value = 58664
$ws.Range("E4").Value = 3050
$ws.Range("F4").Value = 3050
$ws.Range("G4").Value = 2791
$ws.Range("H4").Value = 2101
$ws.Range("I4").Value = 1995
$ws.Range("J4").Value = 106
$ws.Range("K4").Value = 44644
$ws.Range("L4").Value = 29493
$ws.Range("M4").Value = 15151
$ws.Range("N4").Value = 14647
$ws.Range("O4").Value = 504
$ws.Range("P4").Value = 470
$ws.Range("Q4").Value = 3502
$ws.Range("R4").Value = -3662
$ws.Range("S4").Value = 39
$ws.Range("T4").Value = 3365
$ws.Range("U4").Value = 137
$ws.Range("V4").Value = 12895
$ws.Range("W4").Value = 5.2
$ws.Range("X4").Value = 3.58
$ws.Range("Y4").Value = 14.32
$ws.Range("Z4").Value = 4.87
$ws.Range("AA4").Value = 194.67
$ws.Range("AB4").Value = 2674.95
$ws.Range("AC4").Value = 4248
$ws.Range("AD4").Value = 11.04
$ws.Range("AE4").Value = 31280
$ws.Range("AF4").Value = 1.5
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 2.13
$ws.Range("AI4").Value = 23.47
$ws.Range("AJ4").Value = 46957120

# Row 5
$ws.Range("D5").Value = 56847
$ws.Range("E5").Value = 835
$ws.Range("F5").Value = 835
$ws.Range("G5").Value = 193
$ws.Range("H5").Value = 183
$ws.Range("I5").Value = 48
$ws.Range("J5").Value = 135
$ws.Range("K5").Value = 44510
$ws.Range("L5").Value = 30388
$ws.Range("M5").Value = 14123
$ws.Range("N5").Value = 13517
$ws.Range("O5").Value = 605
$ws.Range("P5").Value = 470
$ws.Range("Q5").Value = 2769
$ws.Range("R5").Value = -3774
$ws.Range("S5").Value = 410
$ws.Range("T5").Value = 2857
$ws.Range("U5").Value = -88
$ws.Range("V5").Value = 13391
$ws.Range("W5").Value = 1.47
$ws.Range("X5").Value = 0.32
$ws.Range("Y5").Value = 0.34
$ws.Range("Z5").Value = 0.41
$ws.Range("AA5").Value = 215.17
$ws.Range("AB5").Value = 2568.84
$ws.Range("AC5").Value = 102
$ws.Range("AD5").Value = 605.78
$ws.Range("AE5").Value = 28868
$ws.Range("AF5").Value = 2.14
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 0.32
$ws.Range("AI5").Value = 195.81
$ws.Range("AJ5").Value = 46957120

# Row 6
$ws.Range("D6").Value = 56648
$ws.Range("E6").Value = 1974
$ws.Range("F6").Value = 1974
$ws.Range("G6").Value = 1414
$ws.Range("H6").Value = 1129
$ws.Range("I6").Value = 1057
$ws.Range("K6").Value = 44373
$ws.Range("L6").Value = 29555
$ws.Range("M6").Value = 14818
$ws.Range("N6").Value = 14158
$ws.Range("P6").Value = 470
$ws.Range("Q6").Value = 2869
$ws.Range("R6").Value = -2976
$ws.Range("S6").Value = 100
$ws.Range("T6").Value = 2625
$ws.Range("U6").Value = 244
$ws.Range("V6").Value = 13853
$ws.Range("W6").Value = 3.48
$ws.Range("X6").Value = 1.99
$ws.Range("Y6").Value = 7.64
$ws.Range("Z6").Value = 2.54
$ws.Range("AA6").Value = 199.46
$ws.Range("AB6").Value = 2732.44
$ws.Range("AC6").Value = 2250
$ws.Range("AD6").Value = 12.87
$ws.Range("AE6").Value = 30236
$ws.Range("AF6").Value = 0.96
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 1.73
$ws.Range("AI6").Value = 22.16
$ws.Range("AJ6").Value = 46957120

# Row 7
$ws.Range("D7").Value = 59064
$ws.Range("E7").Value = 2177
$ws.Range("G7").Value = 2058
$ws.Range("H7").Value = 1488
$ws.Range("I7").Value = 1425
$ws.Range("K7").Value = 47356
$ws.Range("L7").Value = 31304
$ws.Range("M7").Value = 16052
$ws.Range("N7").Value = 15479
$ws.Range("P7").Value = 470
$ws.Range("Q7").Value = 3741
$ws.Range("R7").Value = -3384
$ws.Range("S7").Value = 897
$ws.Range("T7").Value = 2493
$ws.Range("U7").Value = 999
$ws.Range("W7").Value = 3.69
$ws.Range("X7").Value = 2.52
$ws.Range("Y7").Value = 9.619999999999999
$ws.Range("Z7").Value = 3.25
$ws.Range("AA7").Value = 195.01
$ws.Range("AC7").Value = 3034
$ws.Range("AD7").Value = 11.06
$ws.Range("AE7").Value = 33058
$ws.Range("AF7").Value = 1.01
$ws.Range("AG7").Value = 654
$ws.Range("AH7").Value = 1.95
$ws.Range("AI7").Value = 21.55

# Row 8
$ws.Range("D8").Value = 63446
$ws.Range("E8").Value = 2630
$ws.Range("G8").Value = 2342
$ws.Range("H8").Value = 1804
$ws.Range("I8").Value = 1735
$ws.Range("K8").Value = 49827
$ws.Range("L8").Value = 32354
$ws.Range("M8").Value = 17474
$ws.Range("N8").Value = 16800
$ws.Range("P8").Value = 470
$ws.Range("Q8").Value = 4223
$ws.Range("R8").Value = -3357
$ws.Range("S8").Value = -290
$ws.Range("T8").Value = 2821
$ws.Range("U8").Value = 2018
$ws.Range("W8").Value = 4.15
$ws.Range("X8").Value = 2.84
$ws.Range("Y8").Value = 10.75
$ws.Range("Z8").Value = 3.71
$ws.Range("AA8").Value = 185.16
$ws.Range("AC8").Value = 3696
$ws.Range("AD8").Value = 9.08
$ws.Range("AE8").Value = 35878
$ws.Range("AF8").Value = 0.9399999999999999
$ws.Range("AG8").Value = 715
$ws.Range("AH8").Value = 2.13
$ws.Range("AI8").Value = 19.36

# Row 9
$ws.Range("D9").Value = 66969
$ws.Range("E9").Value = 2983
$ws.Range("G9").Value = 2650
$ws.Range("H9").Value = 2031
$ws.Range("I9").Value = 1938
$ws.Range("K9").Value = 52518
$ws.Range("L9").Value = 33428
$ws.Range("M9").Value = 19091
$ws.Range("N9").Value = 18298
$ws.Range("P9").Value = 470
$ws.Range("Q9").Value = 4713
$ws.Range("R9").Value = -3532
$ws.Range("S9").Value = -117
$ws.Range("T9").Value = 2913
$ws.Range("U9").Value = 2233
$ws.Range("W9").Value = 4.45
$ws.Range("X9").Value = 3.03
$ws.Range("Y9").Value = 11.04
$ws.Range("Z9").Value = 3.97
$ws.Range("AA9").Value = 175.1
$ws.Range("AC9").Value = 4126
$ws.Range("AD9").Value = 8.130000000000001
$ws.Range("AE9").Value = 39078
$ws.Range("AF9").Value = 0.86
$ws.Range("AG9").Value = 758
$ws.Range("AH9").Value = 2.26
$ws.Range("AI9").Value = 18.36

